$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 693-694, pushing the existing rows 693:783
# down to 695:785 (dimension grows from A1:T783 to A1:T785).
$ws.Rows("693:694").Insert()

# ---- New row 693 ----
$ws.Cells.Item(693, 1).Value2 = 5
$ws.Cells.Item(693, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(693, 3).Value2 = "Maule"
$ws.Cells.Item(693, 4).Value2 = 45131
$ws.Cells.Item(693, 5).Value2 = 7
$ws.Cells.Item(693, 6).Value2 = "Fruta"
$ws.Cells.Item(693, 7).Value2 = 100109
$ws.Cells.Item(693, 8).Value2 = "Uva"
$ws.Cells.Item(693, 9).Value2 = 100109001
$ws.Cells.Item(693, 10).Value2 = "Uva"
$ws.Cells.Item(693, 11).Value2 = "Crimpson Seedless"
$ws.Cells.Item(693, 12).Value2 = "Segunda"
$ws.Cells.Item(693, 13).Value2 = 230
$ws.Cells.Item(693, 14).Value2 = 11000
$ws.Cells.Item(693, 15).Value2 = 11000
$ws.Cells.Item(693, 16).Value2 = 11000
$ws.Cells.Item(693, 17).Value2 = "`$/bandeja 8 kilos"
$ws.Cells.Item(693, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(693, 19).Value2 = 1375
$ws.Cells.Item(693, 20).Value2 = 8

# ---- New row 694 ----
$ws.Cells.Item(694, 1).Value2 = 5
$ws.Cells.Item(694, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(694, 3).Value2 = "Maule"
$ws.Cells.Item(694, 4).Value2 = 45131
$ws.Cells.Item(694, 5).Value2 = 7
$ws.Cells.Item(694, 6).Value2 = "Fruta"
$ws.Cells.Item(694, 7).Value2 = 100109
$ws.Cells.Item(694, 8).Value2 = "Uva"
$ws.Cells.Item(694, 9).Value2 = 100109001
$ws.Cells.Item(694, 10).Value2 = "Uva"
$ws.Cells.Item(694, 11).Value2 = "Red Globe"
$ws.Cells.Item(694, 12).Value2 = "Primera"
$ws.Cells.Item(694, 13).Value2 = 180
$ws.Cells.Item(694, 14).Value2 = 14000
$ws.Cells.Item(694, 15).Value2 = 14000
$ws.Cells.Item(694, 16).Value2 = 14000
$ws.Cells.Item(694, 17).Value2 = "`$/bandeja 8 kilos"
$ws.Cells.Item(694, 18).Value2 = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(694, 19).Value2 = 1750
$ws.Cells.Item(694, 20).Value2 = 8

# Keep column D formatted as a date (style already carried over from the
# row-insert, but set explicitly to be safe).
$ws.Range("D693:D694").NumberFormat = "YYYY-MM-DD HH:MM:SS"
